$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 previously held a bare "Test4" dataset (Test4, j, k, l) with no
# gender/DOB columns. Replace it with a full dataset matching the row 4
# pattern (test_HomePage_FillingForm) but with new name/email/password data.
$ws.Range("A5").Value = "test_HomePage_FillingForm"
$ws.Range("B5").Value = "Jon Doe"
$ws.Range("C5").Value = "jon@test.com"
$ws.Range("D5").Value = "JONDOE123"
$ws.Range("E5").Value = "Male"
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)  # xlPasteFormats: reuse F4's existing date style
$ws.Range("F5").Value2 = $ws.Range("F4").Value2
$excel.CutCopyMode = $false

# Selection moves to the full used range A1:F7 with F7 as the active cell.
$ws.Range("A1:F7").Select()
